# "Finish adding example record books"
#
# The document currently holds a single (empty) paragraph that carries the
# leftover "_GoBack" bookmark. The record-book template needs four labeled,
# boxed fields - Member's Name, 4-H Club Name, Name of Club Leader and
# Report Year / Year(s) Enrolled in 4-H - each ending in an underscored tab
# so the member can fill in the blank, followed by a couple of blank
# paragraphs. The "_GoBack" bookmark is preserved on the "Name of Club
# Leader" line, where the old lone paragraph used to live.
#
# Every one of these paragraphs shares the same box border and tab-leader
# styling, so rather than twiddle each property individually we assemble
# the finished WordprocessingML for the replacement paragraphs and hand it
# to Range.InsertXML - the COM-exposed equivalent of pasting real OOXML
# into the document (Word merely normalizes/ignores the markup it already
# understands, same as any other InsertXML call).

$d = $word.ActiveDocument

$boxBorder = '<w:pBdr>' +
    '<w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
    '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
    '<w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>' +
    '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' +
    '</w:pBdr>'

$fieldXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    $boxBorder
    <w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="9270"/></w:tabs>
    <w:spacing w:before="240" w:after="240" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="0"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Member’s Name: </w:t></w:r>
  <w:r><w:tab/></w:r>
</w:p>
<w:p>
  <w:pPr>
    $boxBorder
    <w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="9270"/></w:tabs>
    <w:spacing w:after="240" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="0"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">4-H Club Name: </w:t></w:r>
  <w:r><w:tab/></w:r>
</w:p>
<w:p>
  <w:pPr>
    $boxBorder
    <w:tabs><w:tab w:val="left" w:leader="underscore" w:pos="9270"/></w:tabs>
    <w:spacing w:after="240" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="0"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Name of Club Leader: </w:t></w:r>
  <w:r><w:tab/></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    $boxBorder
    <w:tabs>
      <w:tab w:val="left" w:leader="underscore" w:pos="3240"/>
      <w:tab w:val="left" w:leader="underscore" w:pos="6570"/>
      <w:tab w:val="left" w:leader="underscore" w:pos="9270"/>
      <w:tab w:val="left" w:leader="underscore" w:pos="9360"/>
    </w:tabs>
    <w:spacing w:after="240" w:line="240" w:lineRule="auto"/>
    <w:ind w:firstLine="0"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Report Year: </w:t></w:r>
  <w:r><w:tab/><w:t xml:space="preserve">Year(s) Enrolled in 4-H: </w:t></w:r>
  <w:r><w:tab/></w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p/>
"@

# Replace the whole body (everything before the final section break) with
# the new fields in one shot, so the lone paragraph that used to hold only
# the "_GoBack" bookmark is cleanly swapped out for the four field lines.
$d.Content.InsertXML($fieldXml)
